$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The player roster (rows 2-15) was re-ordered. Row data (Player, Position, Team)
# stays intact per player, only the row order changes. Rows 16-19 stay as-is.

$newOrder = @(
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Scoot Henderson", "PG", "Portland Trail Blazers"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("D'Angelo Russell", "PG", "Brooklyn Nets"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Norman Powell", "SG,SF", "LA Clippers")
)

for ($i = 0; $i -lt $newOrder.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newOrder[$i][0]
    $ws.Cells.Item($row, 2).Value = $newOrder[$i][1]
    $ws.Cells.Item($row, 3).Value = $newOrder[$i][2]
}
